# fix Sum bugs in Ai model
# Updates the segmentation_results sheet so the revenue / lifetime-value
# "sums" and derived segment/risk labels are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - abdo
$ws.Range("B2").Value = 330
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "High Value"
$ws.Range("F2").Value = 660
$ws.Range("G2").Value = "Low Risk"

# Row 3 - Ahmed
$ws.Range("A3").Value = "Ahmed"
$ws.Range("B3").Value = 1920
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "High Value"
$ws.Range("F3").Value = 5760
$ws.Range("G3").Value = "Low Risk"

# Row 4 - eng.ahmedyaseen4
$ws.Range("A4").Value = "eng.ahmedyaseen4"
$ws.Range("B4").Value = 600
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = "Low Value"
$ws.Range("F4").Value = 600
$ws.Range("G4").Value = "High Risk"

# Row 5 - Mohand
$ws.Range("A5").Value = "Mohand"
$ws.Range("B5").Value = 6000
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = "Low Value"
$ws.Range("F5").Value = 6000
$ws.Range("G5").Value = "High Risk"

# Row 6 - Yaseen
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = "Low Value"
$ws.Range("G6").Value = "High Risk"

# Row 7 - zakria
$ws.Range("B7").Value = 20100
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Medium Value"
$ws.Range("F7").Value = 40200
$ws.Range("G7").Value = "Low Risk"
